# Commit: "test cases for US1570"
# Adds new Surveyor replay-script test data rows to the "Test Environment Data"
# sheet and new survey rows to the "Driver View Test Data" sheet, and leaves
# the workbook's view state (active sheet / selections) the way the author
# left it after editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Users" sheet - the author clicked on B5 while reviewing data.
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B5").Select()

# ---------------------------------------------------------------------------
# 2. "Driver View Test Data" sheet - new survey rows (21-29) for US1570.
# ---------------------------------------------------------------------------
$wsDriver = $wb.Worksheets.Item("Driver View Test Data")

$wsDriver.Range("A22").Value = 21
$wsDriver.Range("B22").Value = "GenerateRandomString(15)"
$wsDriver.Range("C22").Value = "Day"
$wsDriver.Range("D22").Value = "Overcast"
$wsDriver.Range("E22").Value = "Light"
$wsDriver.Range("G22").Value = "Standard"

$wsDriver.Range("A23").Value = 22
$wsDriver.Range("B23").Value = "GenerateRandomString(15)"
$wsDriver.Range("C23").Value = "Night"
$wsDriver.Range("E23").Value = "Light"
$wsDriver.Range("F23").Value = "LessThan50"
$wsDriver.Range("G23").Value = "Standard"

$wsDriver.Range("A24").Value = 23
$wsDriver.Range("B24").Value = "GenerateRandomString(15)"
$wsDriver.Range("C24").Value = "Day"
$wsDriver.Range("D24").Value = "Strong"
$wsDriver.Range("E24").Value = "Light"
$wsDriver.Range("G24").Value = "Standard"

$wsDriver.Range("A25").Value = 24
$wsDriver.Range("B25").Value = "GenerateRandomString(15)"
$wsDriver.Range("C25").Value = "Day"
$wsDriver.Range("D25").Value = "Moderate"
$wsDriver.Range("E25").Value = "Calm"
$wsDriver.Range("G25").Value = "Standard"

$wsDriver.Range("A26").Value = 25
$wsDriver.Range("B26").Value = "GenerateRandomString(15)"
$wsDriver.Range("C26").Value = "Day"
$wsDriver.Range("D26").Value = "Strong"
$wsDriver.Range("E26").Value = "Light"
$wsDriver.Range("G26").Value = "RapidResponse"

$wsDriver.Range("A27").Value = 26
$wsDriver.Range("B27").Value = "GenerateRandomString(15)"
$wsDriver.Range("C27").Value = "Night"
$wsDriver.Range("E27").Value = "Light"
$wsDriver.Range("F27").Value = "GreaterThan50"
$wsDriver.Range("G27").Value = "RapidResponse"

$wsDriver.Range("A28").Value = 27
$wsDriver.Range("B28").Value = "GenerateRandomString(15)"
$wsDriver.Range("C28").Value = "Day"
$wsDriver.Range("D28").Value = "Moderate"
$wsDriver.Range("E28").Value = "Calm"
$wsDriver.Range("G28").Value = "RapidResponse"

$wsDriver.Range("A29").Value = 28
$wsDriver.Range("B29").Value = "GenerateRandomString(15)"
$wsDriver.Range("C29").Value = "Day"
$wsDriver.Range("D29").Value = "Strong"
$wsDriver.Range("E29").Value = "Calm"
$wsDriver.Range("G29").Value = "RapidResponse"

$wsDriver.Range("A30").Value = 29
$wsDriver.Range("B30").Value = "GenerateRandomString(15)"
$wsDriver.Range("C30").Value = "Day"
$wsDriver.Range("D30").Value = "Overcast"
$wsDriver.Range("E30").Value = "Calm"
$wsDriver.Range("G30").Value = "RapidResponse"

$wsDriver.Range("D4").Select()

# ---------------------------------------------------------------------------
# 3. "Test Environment Data" sheet - new Surveyor replay-script rows (10-12).
# ---------------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("Test Environment Data")

$wsEnv.Range("A11").Value = 10
$wsEnv.Range("B11").Value = "SimAuto-Analyzer1"
$wsEnv.Range("C11").Value = "SimAuto-AnalyzerKey1"
$wsEnv.Range("D11").Value = "Surveyor_rr-pic.db3"
$wsEnv.Range("E11").Value = "replay-db3.defn"

$wsEnv.Range("A12").Value = 11
$wsEnv.Range("B12").Value = "SimAuto-Analyzer2"
$wsEnv.Range("C12").Value = "SimAuto-AnalyzerKey2"
$wsEnv.Range("D12").Value = "Surveyor_rr-sqacudr.db3"
$wsEnv.Range("E12").Value = "replay-db3.defn"

$wsEnv.Range("A13").Value = 12
$wsEnv.Range("B13").Value = "SimAuto-Analyzer1"
$wsEnv.Range("C13").Value = "SimAuto-AnalyzerKey1"
$wsEnv.Range("D13").Value = "Surveyor_rr.db3"
$wsEnv.Range("E13").Value = "replay-db3.defn"

# This is the last sheet the author interacted with, so it ends up the
# active tab of the workbook.
$wsEnv.Range("A13").Select()
